$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("latest_eval")

$ws.Range("B2").Value = 0.001974836808913309
$ws.Range("C2").Value = 0.7564986997685976
$ws.Range("D2").Value = 1.331115891817462
$ws.Range("E2").Value = 1.153739958490414
$ws.Range("F2").Value = 1.165451647513039
$ws.Range("G2").Value = 50
$ws.Range("B3").Value = 0.1715060524283095
$ws.Range("C3").Value = 0.7455020267087988
$ws.Range("D3").Value = 1.556349277421864
$ws.Range("E3").Value = 1.247537284982643
$ws.Range("F3").Value = 1.248768565063892
$ws.Range("G3").Value = 48
$ws.Range("B4").Value = -0.1210599245473831
$ws.Range("C4").Value = 0.6932517703744887
$ws.Range("D4").Value = 1.42325916890264
$ws.Range("E4").Value = 1.19300426189626
$ws.Range("F4").Value = 1.199405667105797
$ws.Range("G4").Value = 48
$ws.Range("B5").Value = 0.1732007841850449
$ws.Range("C5").Value = 0.7115279881092994
$ws.Range("D5").Value = 1.411811716246263
$ws.Range("E5").Value = 1.188196833965763
$ws.Range("F5").Value = 1.188494915725877
$ws.Range("G5").Value = 46
$ws.Range("B6").Value = 0.03614660433360007
$ws.Range("C6").Value = 0.5611928274455102
$ws.Range("D6").Value = 0.7065544791159908
$ws.Range("E6").Value = 0.8405679503264389
$ws.Range("F6").Value = 0.849279849945092
$ws.Range("G6").Value = 45
$ws.Range("B7").Value = -0.03741410742334574
$ws.Range("C7").Value = 0.6726782106122756
$ws.Range("D7").Value = 1.306474751485271
$ws.Range("E7").Value = 1.143011264811188
$ws.Range("F7").Value = 1.159578656834397
$ws.Range("G7").Value = 34
$ws.Range("B8").Value = -0.04498115623876454
$ws.Range("C8").Value = 0.6806687468339884
$ws.Range("D8").Value = 1.217036297618114
$ws.Range("E8").Value = 1.103193680918321
$ws.Range("F8").Value = 1.119366851947603
$ws.Range("G8").Value = 33
$ws.Range("B9").Value = 0.008498424339290329
$ws.Range("C9").Value = 0.6040240707239352
$ws.Range("D9").Value = 0.6504015573506643
$ws.Range("E9").Value = 0.8064747716765009
$ws.Range("F9").Value = 0.832877315741465
$ws.Range("G9").Value = 16
$ws.Range("B10").Value = -0.005607768909606597
$ws.Range("C10").Value = 0.504936332959623
$ws.Range("D10").Value = 0.4301487342848143
$ws.Range("E10").Value = 0.6558572514540144
$ws.Range("F10").Value = 0.6913089735443678
$ws.Range("G10").Value = 10